# Append two new rows (130 and 131) of feed log data to Sheet1, mirroring
# the existing data rows (run_id, rss_url_id, date, response, item_count).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 130
$ws.Cells.Item(130, 1).Value = 129
$ws.Cells.Item(130, 2).Value = 1
$ws.Cells.Item(130, 3).Value = "2024-06-17 16:16:16"
$ws.Cells.Item(130, 4).Value = 200
$ws.Cells.Item(130, 5).Value = 16

# Row 131
$ws.Cells.Item(131, 1).Value = 130
$ws.Cells.Item(131, 2).Value = 2
$ws.Cells.Item(131, 3).Value = "2024-06-17 16:16:16"
$ws.Cells.Item(131, 4).Value = 200
$ws.Cells.Item(131, 5).Value = 1
